$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("borehole")

# Rename the "temperature_accuracy" column header to "temperature_uncertainty"
$ws.Range("O1").Value = "temperature_uncertainty"

# Update the column's comment to describe the renamed/generalized field
$newComment = "[number: °C] temperature_uncertainty`n`nEstimated temperature uncertainty (as reported)."
$ws.Range("O1").Comment.Text($newComment) | Out-Null

# Widen column O slightly to fit the longer header/description
$ws.Columns.Item(15).ColumnWidth = 27.5
